$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.357105851173401
$ws.Range("B1").Value = 1.230889797210693
$ws.Range("C1").Value = 3.218952417373657
$ws.Range("D1").Value = 2.875737190246582
$ws.Range("E1").Value = 0.907846987247467
